{"js": "// Split the single Bibliografia reference run into 15 runs, one per\n// numbered citation, joined by line breaks (<w:br/>), matching the\n// source diff. Each citation keeps its own <w:t> (with a trailing\n// space preserved via xml:space=\"preserve\" where present in the\n// original text), separated from the next by a break - all inside a\n// single run, exactly as in the target OOXML.\n\nconst REFERENCES = [\"1.CORR\u00caA, H. L.; GIANESI, I. G. N.; CAON, M. Planejamento, programa\u00e7\u00e3o e controle da produ\u00e7\u00e3o: MRPII/ERP conceitos, uso e implanta\u00e7\u00e3o. 5. ed. S\u00e3o Paulo: Atlas, 2007. \", \"2.CORR\u00caA, H. L.; CORR\u00caA, C. A. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o e Opera\u00e7\u00f5es: manufatura e servi\u00e7os: uma abordagem estrat\u00e9gica. 2.ed. S\u00e3o Paulo: Atlas, 2011. \", \"3.DAVIS, M.M. et al. Fundamentos da administra\u00e7\u00e3o da Produ\u00e7\u00e3o. Porto Alegre: Bookman, 2018. \", \"4.FERNANDES, F.C.F.; GODINHO FILHO. Planejamento e controle da produ\u00e7\u00e3o: dos fundamentos ao essencial. S\u00e3o Paulo: Atlas, 2010. \", \"5.GAITHER, N.; FRAZIER, G. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o e Opera\u00e7\u00f5es. 8. ed. S\u00e3o Paulo: Pioneira Thomson, 2005. \", \"6.GON\u00c7ALVES, P.S. Administra\u00e7\u00e3o de materiais. Rio de Janeiro: Elsevier, 2013. \", \"7.HEIZER, J.; RENDER, B. Administra\u00e7\u00e3o de Opera\u00e7\u00f5es: bens e servi\u00e7os. 5. ed. Rio de Janeiro: LTC, 2001.\", \"8.JACOBS, F.R.; CHASE, R. B. Administra\u00e7\u00e3o da produ\u00e7\u00e3o e de opera\u00e7\u00f5es: o essencial. Porto Alegre: Bookman, 2009. \", \"9.LUSTOSA, L. et. Al. Planejamento e controle da produ\u00e7\u00e3o. Rio de Janeiro: Elsevier, 2008. \", \"10.MOREIRA, D. A. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o e Opera\u00e7\u00f5es. 2. ed. S\u00e3o Paulo: Cengage Learning, 2008.\", \"11.REID, R.D.; SANDERS, N. R. Gest\u00e3o de opera\u00e7\u00f5es. Rio de Janeiro: LTC, 2005. \", \"12.SLACK, N., BRANDON-JONES, A., JOHNSTON, R. Administra\u00e7\u00e3o da produ\u00e7\u00e3o. Henrique Luiz Corr\u00eaa (Trad.). 3. ed. S\u00e3o Paulo: Atlas, 2018. \", \"13.TUBINO, D.V. Planejamento e controle da produ\u00e7\u00e3o: teoria e pr\u00e1tica. 2.ed. S\u00e3o Paulo: Atlas, 2009. \", \"14.VOLLMANN, T.; BERRY, W.; WHYBARK, D.; JACOBS, F. Sistemas de planejamento e controle da produ\u00e7\u00e3o: para o gerenciamento da cadeia de suprimentos. 5. ed. Porto Alegre: Bookman, 2006. \", \"15.WANKE, P.F. Ger\u00eancia de opera\u00e7\u00f5es: uma abordagem log\u00edstica. S\u00e3o Paulo: Atlas, 2010.\"];\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Build the flat-OPC WordprocessingML package insertOoxml() expects:\n// a <pkg:package> wrapping a single /word/document.xml part whose body\n// holds exactly the replacement paragraph.\nfunction buildFlatOpcParagraph(items) {\n  const runInner = items\n    .map((text, i) => {\n      const preserve = /^\\s|\\s$/.test(text);\n      const attr = preserve ? ' xml:space=\"preserve\"' : \"\";\n      const t = `<w:t${attr}>${escapeXml(text)}</w:t>`;\n      return i < items.length - 1 ? `${t}<w:br/>` : t;\n    })\n    .join(\"\");\n\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    `<w:body><w:p><w:r>${runInner}</w:r></w:p></w:body>` +\n    \"</w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\n// Locate the Bibliografia reference paragraph by its distinctive leading\n// text rather than a hard-coded index.\nconst marker = REFERENCES[0].slice(0, 20);\nconst target = paragraphs.items.find((p) => p.text.startsWith(marker));\nif (!target) {\n  throw new Error(\"Could not find the Bibliografia references paragraph\");\n}\n\nconst ooxml = buildFlatOpcParagraph(REFERENCES);\ntarget.getRange(\"Whole\").insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n\n", "ps1": "$d = $word.ActiveDocument\n\n# Each numbered citation from the Bibliografia paragraph, to be joined\n# back together with <w:br/> line breaks inside a single run (matching\n# the source diff), instead of the original one-run, no-breaks text.\n$references = @(\n    '1.CORR\u00caA, H. L.; GIANESI, I. G. N.; CAON, M. Planejamento, programa\u00e7\u00e3o e controle da produ\u00e7\u00e3o: MRPII/ERP conceitos, uso e implanta\u00e7\u00e3o. 5. ed. S\u00e3o Paulo: Atlas, 2007. ',\n    '2.CORR\u00caA, H. L.; CORR\u00caA, C. A. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o e Opera\u00e7\u00f5es: manufatura e servi\u00e7os: uma abordagem estrat\u00e9gica. 2.ed. S\u00e3o Paulo: Atlas, 2011. ',\n    '3.DAVIS, M.M. et al. Fundamentos da administra\u00e7\u00e3o da Produ\u00e7\u00e3o. Porto Alegre: Bookman, 2018. ',\n    '4.FERNANDES, F.C.F.; GODINHO FILHO. Planejamento e controle da produ\u00e7\u00e3o: dos fundamentos ao essencial. S\u00e3o Paulo: Atlas, 2010. ',\n    '5.GAITHER, N.; FRAZIER, G. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o e Opera\u00e7\u00f5es. 8. ed. S\u00e3o Paulo: Pioneira Thomson, 2005. ',\n    '6.GON\u00c7ALVES, P.S. Administra\u00e7\u00e3o de materiais. Rio de Janeiro: Elsevier, 2013. ',\n    '7.HEIZER, J.; RENDER, B. Administra\u00e7\u00e3o de Opera\u00e7\u00f5es: bens e servi\u00e7os. 5. ed. Rio de Janeiro: LTC, 2001.',\n    '8.JACOBS, F.R.; CHASE, R. B. Administra\u00e7\u00e3o da produ\u00e7\u00e3o e de opera\u00e7\u00f5es: o essencial. Porto Alegre: Bookman, 2009. ',\n    '9.LUSTOSA, L. et. Al. Planejamento e controle da produ\u00e7\u00e3o. Rio de Janeiro: Elsevier, 2008. ',\n    '10.MOREIRA, D. A. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o e Opera\u00e7\u00f5es. 2. ed. S\u00e3o Paulo: Cengage Learning, 2008.',\n    '11.REID, R.D.; SANDERS, N. R. Gest\u00e3o de opera\u00e7\u00f5es. Rio de Janeiro: LTC, 2005. ',\n    '12.SLACK, N., BRANDON-JONES, A., JOHNSTON, R. Administra\u00e7\u00e3o da produ\u00e7\u00e3o. Henrique Luiz Corr\u00eaa (Trad.). 3. ed. S\u00e3o Paulo: Atlas, 2018. ',\n    '13.TUBINO, D.V. Planejamento e controle da produ\u00e7\u00e3o: teoria e pr\u00e1tica. 2.ed. S\u00e3o Paulo: Atlas, 2009. ',\n    '14.VOLLMANN, T.; BERRY, W.; WHYBARK, D.; JACOBS, F. Sistemas de planejamento e controle da produ\u00e7\u00e3o: para o gerenciamento da cadeia de suprimentos. 5. ed. Porto Alegre: Bookman, 2006. ',\n    '15.WANKE, P.F. Ger\u00eancia de opera\u00e7\u00f5es: uma abordagem log\u00edstica. S\u00e3o Paulo: Atlas, 2010.'\n)\n\nfunction Convert-XmlText {\n    param([string]$Text)\n    return $Text.Replace('&', '&amp;').Replace('<', '&lt;').Replace('>', '&gt;')\n}\n\n$runInner = \"\"\nfor ($i = 0; $i -lt $references.Count; $i++) {\n    $text = $references[$i]\n    $escaped = Convert-XmlText $text\n    if ($text -match '^\\s' -or $text -match '\\s$') {\n        $runInner += \"<w:t xml:space=`\"preserve`\">$escaped</w:t>\"\n    } else {\n        $runInner += \"<w:t>$escaped</w:t>\"\n    }\n    if ($i -lt $references.Count - 1) {\n        $runInner += \"<w:br/>\"\n    }\n}\n\n$flatOpc = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + \n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + \n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + \n    \"<w:body><w:p><w:r>$runInner</w:r></w:p></w:body>\" + \n    '</w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Locate the Bibliografia reference paragraph (identified by its leading\n# \"1.CORR\u00caA\" marker) instead of a hard-coded index, so the script is\n# resilient to minor structural differences elsewhere in the document.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith($references[0].Substring(0, 20))) {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Bibliografia paragraph not found\"\n}\n\n[void]$target.Range.InsertXML($flatOpc)\n"}
